$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 299, shifting existing rows 299-337 down to 300-338.
$ws.Rows.Item(299).Insert()

# Populate the newly inserted row 299 with the new record.
$ws.Range("A299").Value = 10
$ws.Range("B299").Value = "Vega Modelo de Temuco"
$ws.Range("C299").Value = "La Araucanía"
$ws.Range("D299").Value = 44578
$ws.Range("E299").Value = 9
$ws.Range("F299").Value = 100112032
$ws.Range("G299").Value = "Zapallo italiano"
$ws.Range("H299").Value = "Sin especificar"
$ws.Range("I299").Value = "Primera"
$ws.Range("J299").Value = 300
$ws.Range("K299").Value = 12000
$ws.Range("L299").Value = 12000
$ws.Range("M299").Value = 12000
$ws.Range("N299").Value = "$/caja 60 unidades"
$ws.Range("O299").Value = "Región de Arica y Parinacota"
$ws.Range("P299").Value = 200
$ws.Range("Q299").Value = 60
$ws.Range("R299").Value = "Hortaliza"
